# Updates the "Price" (D) and "Volume(1h)" (E) columns for each coin row with
# freshly scraped figures, and reorders two coin pairs (Dai/Avalanche at rows
# 20-21, Aptos/WEMIXTOKEN at rows 41-42) while keeping the A-column rank in
# place. All D/E cells in the sheet are stored as text (not numbers), so every
# Price value is written with a leading apostrophe text-prefix to stop Excel
# from re-parsing strings like "3.670" or "0.00001056" as numbers (which would
# silently drop the trailing zero / switch to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.238.00'
$ws.Range('E2').Value = '  -0.15%  '

$ws.Range('D3').Value = '''1.787.89'
$ws.Range('E3').Value = '  -0.76%  '

$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').Value = '''316.24'
$ws.Range('E5').Value = '  -0.05%  '

$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('D7').Value = '''0.5316'
$ws.Range('E7').Value = '  -2.92%  '

$ws.Range('D8').Value = '''0.3753'
$ws.Range('E8').Value = '  -2.55%  '

$ws.Range('D9').Value = '''0.07481'
$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').Value = '''41.50'
$ws.Range('E10').Value = '  -3.26%  '

$ws.Range('D11').Value = '''1.092'
$ws.Range('E11').Value = '  -2.52%  '

$ws.Range('D12').Value = '''1.003'
$ws.Range('E12').Value = '  +0.23%  '

$ws.Range('E13').Value = '  -3.35%  '

$ws.Range('D14').Value = '''6.093'
$ws.Range('E14').Value = '  -1.68%  '

$ws.Range('D15').Value = '''7.237'
$ws.Range('E15').Value = '  -1.34%  '

$ws.Range('D16').Value = '''1.766.51'
$ws.Range('E16').Value = '  -1.68%  '

$ws.Range('D17').Value = '''89.14'
$ws.Range('E17').Value = '  -2.78%  '

$ws.Range('D18').Value = '''0.00001056'
$ws.Range('E18').Value = '  -1.19%  '

$ws.Range('D19').Value = '''0.06493'
$ws.Range('E19').Value = '  +0.71%  '

$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.001'
$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '''17.41'
$ws.Range('E21').Value = '  +0.99%  '

$ws.Range('D22').Value = '''5.924'
$ws.Range('E22').Value = '  -0.77%  '

$ws.Range('D23').Value = '''28.268.86'
$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('D24').Value = '''11.11'
$ws.Range('E24').Value = '  -3.16%  '

$ws.Range('D25').Value = '''2.092'
$ws.Range('E25').Value = '  -3.62%  '

$ws.Range('D26').Value = '''158.12'
$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('D27').Value = '''20.27'
$ws.Range('E27').Value = '  -1.74%  '

$ws.Range('D28').Value = '''1.975.96'
$ws.Range('E28').Value = '  -1.51%  '

$ws.Range('D29').Value = '''2.288'
$ws.Range('E29').Value = '  -5.49%  '

$ws.Range('D30').Value = '''121.03'
$ws.Range('E30').Value = '  -2.06%  '

$ws.Range('D31').Value = '''1.096'
$ws.Range('E31').Value = '  -3.87%  '

$ws.Range('D32').Value = '''0.1043'
$ws.Range('E32').Value = '  +2.75%  '

$ws.Range('D33').Value = '''3.661'
$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('E34').Value = '  -3.64%  '

$ws.Range('D35').Value = '''0.2243'
$ws.Range('E35').Value = '  -2.91%  '

$ws.Range('D36').Value = '''0.06392'
$ws.Range('E36').Value = '  +1.28%  '

$ws.Range('D37').Value = '''0.02273'
$ws.Range('E37').Value = '  -1.91%  '

$ws.Range('D38').Value = '''4.992'
$ws.Range('E38').Value = '  -0.56%  '

$ws.Range('D39').Value = '''8.468'
$ws.Range('E39').Value = '  -5.03%  '

$ws.Range('D40').Value = '''0.6167'
$ws.Range('E40').Value = '  -3.25%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '''11.02'
$ws.Range('E41').Value = '  -4.66%  '

$ws.Range('B42').Value = 'WEMIXTOKEN'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '''1.430'
$ws.Range('E42').Value = '  +3.40%  '

$ws.Range('D43').Value = '''1.176'
$ws.Range('E43').Value = '  +1.34%  '

$ws.Range('D44').Value = '''1.002'
$ws.Range('E44').Value = '  +0.19%  '

$ws.Range('D45').Value = '''13.27'
$ws.Range('E45').Value = '  -1.87%  '

$ws.Range('D46').Value = '''3.670'
$ws.Range('E46').Value = '  -0.17%  '

$ws.Range('E47').Value = '  -3.58%  '

$ws.Range('D48').Value = '''125.73'
$ws.Range('E48').Value = '  +1.45%  '

$ws.Range('D49').Value = '''1.200'
$ws.Range('E49').Value = '  +4.74%  '

$ws.Range('D50').Value = '''1.929'
$ws.Range('E50').Value = '  -2.16%  '

$ws.Range('D51').Value = '''0.06840'
$ws.Range('E51').Value = '  -0.76%  '
